$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.606.16"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.043.41"
$ws.Range("E3").Value = "  -2.94%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.09"
$ws.Range("E5").Value = "  -4.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "602.79"
$ws.Range("E6").Value = "  -2.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.08"
$ws.Range("E7").Value = "  -1.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.371"
$ws.Range("E8").Value = "  -9.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.791"
$ws.Range("E10").Value = "  +6.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.039.83"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.194"
$ws.Range("E12").Value = "  -4.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.446.61"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000236"
$ws.Range("E14").Value = "  -8.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.16"
$ws.Range("E15").Value = "  -5.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("E16").Value = "  -4.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.614.33"
$ws.Range("E17").Value = "  -2.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.056.02"
$ws.Range("E18").Value = "  -2.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.49"
$ws.Range("E19").Value = "  -8.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.15"
$ws.Range("E20").Value = "  -5.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.61"
$ws.Range("E21").Value = "  -4.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.62"
$ws.Range("E22").Value = "  -4.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.69"
$ws.Range("E23").Value = "  -8.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000186"
$ws.Range("E24").Value = "  -11.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.33"
$ws.Range("E25").Value = "  +5.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.42"
$ws.Range("E26").Value = "  -8.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "84.26"
$ws.Range("E27").Value = "  -4.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.57"
$ws.Range("E28").Value = "  -3.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.220.11"
$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("E31").Value = "  +8.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.175"
$ws.Range("E32").Value = "  +2.90%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.10"
$ws.Range("E33").Value = "  +10.55%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.121"
$ws.Range("E34").Value = "  -11.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.92"
$ws.Range("E35").Value = "  -4.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.53"
$ws.Range("E36").Value = "  -8.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.153"
$ws.Range("E37").Value = "  -6.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.98"
$ws.Range("E38").Value = "  -5.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  -2.14%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.93"
$ws.Range("E40").Value = "  +3.53%  "

$ws.Range("B41").Value = "MantraDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.78"
$ws.Range("E41").Value = "  -4.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.433"
$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "455.24"
$ws.Range("E43").Value = "  -6.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.22"
$ws.Range("E44").Value = "  -7.15%  "

$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.07"
$ws.Range("E46").Value = "  -11.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.45"
$ws.Range("E47").Value = "  -1.25%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.660"
$ws.Range("E48").Value = "  -5.73%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.80"
$ws.Range("E49").Value = "  -8.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.67"
$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  +0.19%  "
